# Sensors.xlsx update: "Added support for AM2320 sensor by I2C. Optimized work with I2C"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 5 (DHT12): Protocol "Single wire, I2C" (I2C highlighted red) -> "Single wire" ---
$ws.Range("C5").Value = "Single wire"

# --- Row 6 (DHT21 / AM2301): power range, accuracy and resolution corrections ---
$ws.Range("B6").Value = "3.3 - 5 V"

# Accuracy becomes " ±1 °C" (plain/default black, was red-tinted " ±0.5 °C")
$ws.Range("E6").Value = " ±1 °C"
$accChars = $ws.Range("E6").Characters(1, 4)
$accChars.Font.ColorIndex = -4105

$ws.Range("I6").Value = "0.1 %"

# --- Row 8 (AM2320, new sensor): Protocol keeps "Single wire, I2C" text but drop red I2C highlight ---
$ws.Range("C8").Value = "Single wire, I2C"
$protoChars = $ws.Range("C8").Characters(14, 3)
$protoChars.Font.ColorIndex = -4105

# --- Row 10 (LM75): power range + resolution corrections ---
$ws.Range("B10").Value = "2.8 - 5.5 V"
$ws.Range("F10").Value = "0.1 °C"

# --- Row 11 (BMP280): Protocol "I2C, SPI" (SPI highlighted red) -> "I2C" (SPI support dropped) ---
$ws.Range("C11").Value = "I2C"
